# Auto-generated edit script: updates LevePrice/LeveProfit derived columns (H-N)
# across the 8 crafting-job worksheets, per the scheduled price-data refresh.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (item id 5489)
$ws.Range("H2").Value = 4565.95
$ws.Range("I2").Value = 2755
$ws.Range("J2").Value = 9998.799999999999
$ws.Range("K2").Value = 2755
$ws.Range("L2").Value = 9998.799999999999
$ws.Range("M2").Value = -2642
$ws.Range("N2").Value = -10224.8
# Row 62 (item id 27781)
$ws.Range("H62").Value = 4189.2
$ws.Range("I62").Value = 4189.2
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4189.2
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3565.2
$ws.Range("N62").ClearContents()
# Row 65 (item id 27781)
$ws.Range("H65").Value = 4189.2
$ws.Range("I65").Value = 4189.2
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20946
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -17826
$ws.Range("N65").ClearContents()
# Row 80 (item id 12605)
$ws.Range("H80").Value = 20833972
$ws.Range("J80").Value = 22727906
$ws.Range("L80").Value = 68183718
$ws.Range("N80").Value = -68185714
# Row 83 (item id 12605)
$ws.Range("H83").Value = 20833972
$ws.Range("J83").Value = 22727906
$ws.Range("L83").Value = 204551154
$ws.Range("N83").Value = -204561138
# Row 98 (item id 36237)
$ws.Range("H98").Value = 3600.8125
$ws.Range("I98").Value = 2832.889
$ws.Range("K98").Value = 2832.889
$ws.Range("M98").Value = -1334.889
# Row 105 (item id 18668)
$ws.Range("H105").Value = 36500
$ws.Range("J105").Value = 36500
$ws.Range("L105").Value = 36500
$ws.Range("N105").Value = -43488
# Row 107 (item id 27766)
$ws.Range("H107").Value = 2703.9375
$ws.Range("I107").Value = 2178.0908
$ws.Range("J107").Value = 3860.8
$ws.Range("K107").Value = 2178.0908
$ws.Range("L107").Value = 3860.8
$ws.Range("M107").Value = -258.0907999999999
$ws.Range("N107").Value = -7700.8
# Row 113 (item id 27775)
$ws.Range("H113").Value = 2716.5
$ws.Range("I113").Value = 1659.8
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 1659.8
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = 1594.2
$ws.Range("N113").Value = -14508
# Row 122 (item id 36237)
$ws.Range("H122").Value = 3600.8125
$ws.Range("I122").Value = 2832.889
$ws.Range("K122").Value = 8498.667000000001
$ws.Range("M122").Value = -6048.667000000001
# Row 137 (item id 44013)
$ws.Range("H137").Value = 1339406.1
$ws.Range("I137").Value = 4311.552
$ws.Range("K137").Value = 12934.656
$ws.Range("M137").Value = -10384.656

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5 (item id 5091)
$ws.Range("H5").Value = 728.25
$ws.Range("I5").Value = 763.2857
$ws.Range("J5").Value = 701
$ws.Range("K5").Value = 763.2857
$ws.Range("L5").Value = 701
$ws.Range("M5").Value = -651.2857
$ws.Range("N5").Value = -925
# Row 61 (item id 43999)
$ws.Range("H61").Value = 815933.25
$ws.Range("I61").Value = 20725.291
$ws.Range("J61").Value = 4608463.5
$ws.Range("K61").Value = 20725.291
$ws.Range("L61").Value = 4608463.5
$ws.Range("M61").Value = -20513.291
$ws.Range("N61").Value = -4608887.5
# Row 122 (item id 36168)
$ws.Range("H122").Value = 2355.5833
$ws.Range("I122").Value = 2355.5833
$ws.Range("K122").Value = 7066.749899999999
$ws.Range("M122").Value = -4616.749899999999
# Row 136 (item id 43999)
$ws.Range("H136").Value = 815933.25
$ws.Range("I136").Value = 20725.291
$ws.Range("J136").Value = 4608463.5
$ws.Range("K136").Value = 62175.87300000001
$ws.Range("L136").Value = 13825390.5
$ws.Range("M136").Value = -59625.87300000001
$ws.Range("N136").Value = -13830490.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4 (item id 5091)
$ws.Range("H4").Value = 728.25
$ws.Range("I4").Value = 763.2857
$ws.Range("J4").Value = 701
$ws.Range("K4").Value = 763.2857
$ws.Range("L4").Value = 701
$ws.Range("M4").Value = -648.2857
$ws.Range("N4").Value = -931
# Row 99 (item id 19943)
$ws.Range("H99").Value = 8813.625
$ws.Range("I99").Value = 9321.200000000001
$ws.Range("K99").Value = 9321.200000000001
$ws.Range("M99").Value = -7823.200000000001
# Row 105 (item id 19947)
$ws.Range("H105").Value = 12444.667
$ws.Range("I105").Value = 9157.053
$ws.Range("J105").Value = 20252.75
$ws.Range("K105").Value = 9157.053
$ws.Range("L105").Value = 20252.75
$ws.Range("M105").Value = -7410.053
$ws.Range("N105").Value = -23746.75
# Row 134 (item id 43998)
$ws.Range("H134").Value = 27274484
$ws.Range("J134").Value = 225002620
$ws.Range("L134").Value = 675007860
$ws.Range("N134").Value = -675012930

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (item id 27691)
$ws.Range("H16").Value = 5999
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5999
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5999
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -6573
# Row 31 (item id 44023)
$ws.Range("H31").Value = 2585.3555
$ws.Range("I31").Value = 2838.6956
$ws.Range("J31").Value = 2320.5
$ws.Range("K31").Value = 2838.6956
$ws.Range("L31").Value = 2320.5
$ws.Range("M31").Value = -2543.6956
$ws.Range("N31").Value = -2910.5
# Row 34 (item id 44023)
$ws.Range("H34").Value = 2585.3555
$ws.Range("I34").Value = 2838.6956
$ws.Range("J34").Value = 2320.5
$ws.Range("K34").Value = 2838.6956
$ws.Range("L34").Value = 2320.5
$ws.Range("M34").Value = -2636.6956
$ws.Range("N34").Value = -2724.5
# Row 105 (item id 19928)
$ws.Range("H105").Value = 2742
$ws.Range("I105").Value = 1675.5
$ws.Range("K105").Value = 1675.5
$ws.Range("M105").Value = 71.5
# Row 113 (item id 27691)
$ws.Range("H113").Value = 5999
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5999
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5999
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10339
# Row 134 (item id 44020)
$ws.Range("H134").Value = 2404.6562
$ws.Range("I134").Value = 2064.524
$ws.Range("J134").Value = 3054
$ws.Range("K134").Value = 6193.572
$ws.Range("L134").Value = 9162
$ws.Range("M134").Value = -3658.572
$ws.Range("N134").Value = -14232
# Row 139 (item id 43258)
$ws.Range("H139").Value = 55889
$ws.Range("I139").Value = 55889
$ws.Range("K139").Value = 55889
$ws.Range("M139").Value = -50749

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 103 (item id 19839)
$ws.Range("H103").Value = 591355.9399999999
$ws.Range("I103").Value = 1428785.2
$ws.Range("J103").Value = 5155.4
$ws.Range("K103").Value = 4286355.6
$ws.Range("L103").Value = 15466.2
$ws.Range("M103").Value = -4285476.6
$ws.Range("N103").Value = -17224.2
# Row 106 (item id 19819)
$ws.Range("H106").Value = 9150.429
$ws.Range("J106").Value = 9150.429
$ws.Range("L106").Value = 27451.287
$ws.Range("N106").Value = -29343.287
# Row 115 (item id 27861)
$ws.Range("H115").Value = 964
$ws.Range("I115").Value = 964
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2892
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1717
$ws.Range("N115").ClearContents()
# Row 132 (item id 43972)
$ws.Range("H132").Value = 1569
$ws.Range("I132").Value = 1051.8889
$ws.Range("K132").Value = 9467.000099999999
$ws.Range("M132").Value = -6937.000099999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (item id 14146)
$ws.Range("H70").Value = 4945.636
$ws.Range("I70").Value = 4599.75
$ws.Range("K70").Value = 4599.75
$ws.Range("M70").Value = -4329.75
# Row 73 (item id 14146)
$ws.Range("H73").Value = 4945.636
$ws.Range("I73").Value = 4599.75
$ws.Range("K73").Value = 4599.75
$ws.Range("M73").Value = -3663.75
# Row 107 (item id 27802)
$ws.Range("H107").Value = 84045
$ws.Range("J107").Value = 1124.2858
$ws.Range("L107").Value = 1124.2858
$ws.Range("N107").Value = -4964.2858
# Row 123 (item id 34150)
$ws.Range("H123").Value = 47999.668
$ws.Range("J123").Value = 51999.5
$ws.Range("L123").Value = 51999.5
$ws.Range("N123").Value = -56899.5
# Row 126 (item id 36184)
$ws.Range("H126").Value = 2584.4443
$ws.Range("J126").Value = 2442
$ws.Range("L126").Value = 7326
$ws.Range("N126").Value = -12266
# Row 132 (item id 44008)
$ws.Range("H132").Value = 755934.0600000001
$ws.Range("I132").Value = 1141.875
$ws.Range("J132").Value = 988177.9
$ws.Range("K132").Value = 3425.625
$ws.Range("L132").Value = 2964533.7
$ws.Range("M132").Value = -895.625
$ws.Range("N132").Value = -2969593.7

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55 (item id 5284)
$ws.Range("H55").Value = 588.9143
$ws.Range("J55").Value = 773.4761999999999
$ws.Range("L55").Value = 773.4761999999999
$ws.Range("N55").Value = -1119.4762
# Row 100 (item id 19995)
$ws.Range("H100").Value = 4164.44
$ws.Range("I100").Value = 3290.5557
$ws.Range("J100").Value = 4656
$ws.Range("K100").Value = 3290.5557
$ws.Range("L100").Value = 4656
$ws.Range("M100").Value = -2749.5557
$ws.Range("N100").Value = -5738

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132 (item id 44029)
$ws.Range("H132").Value = 25001892
$ws.Range("I132").Value = 35715584
$ws.Range("J132").Value = 3273.6667
$ws.Range("K132").Value = 107146752
$ws.Range("L132").Value = 9821.000100000001
$ws.Range("M132").Value = -107144222
$ws.Range("N132").Value = -14881.0001

